$wb = $excel.ActiveWorkbook

# Switch to the "choices" sheet and update its header cell C1
# from "display.text" to "display.title.text".
$choices = $wb.Worksheets.Item("choices")
$choices.Activate()
$choices.Range("C1").Value = "display.title.text"
$choices.Range("C2").Select()
